# Weekly fruit/vegetable price update:
# Insert one new daily-price record for "Jengibre" (ginger) at row 74,
# pushing the existing rows 74-109 down to 75-110 (dimension grows to A1:R110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 74 (shifts 74..109 -> 75..110).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new record.
$ws.Cells.Item(74, 1).Value  = 6
$ws.Cells.Item(74, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(74, 3).Value  = "Metropolitana"
$ws.Cells.Item(74, 4).Value  = 44873
$ws.Cells.Item(74, 5).Value  = 13
$ws.Cells.Item(74, 6).Value  = 100114007
$ws.Cells.Item(74, 7).Value  = "Jengibre"
$ws.Cells.Item(74, 8).Value  = "Sin especificar"
$ws.Cells.Item(74, 9).Value  = "Primera"
$ws.Cells.Item(74, 10).Value = 280
$ws.Cells.Item(74, 11).Value = 12000
$ws.Cells.Item(74, 12).Value = 13000
$ws.Cells.Item(74, 13).Value = 12536
$ws.Cells.Item(74, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(74, 15).Value = "Perú"
$ws.Cells.Item(74, 16).Value = 964
$ws.Cells.Item(74, 17).Value = 13
$ws.Cells.Item(74, 18).Value = "Hortaliza"
